# Generate Report for Handoff
#
# This mirrors a re-run of the localization-status report generator: a few
# "Latest Handoff/Handback Datetime" timestamps that used to be distinct
# values have now collapsed onto a later, shared timestamp (the handoff/
# handback apparently completed a bit later than before). Concretely:
#   - Overview sheet, column D ("Latest Handoff Date"):
#       "2016-24-20 20:24:12" -> "2016-24-20 20:24:57"
#       (this also absorbs rows that used to read "2016-24-20 20:24:35")
#   - zh-cn sheet, column E ("Latest Handoff Datetime"):
#       "2016-03-20 20:24:08" -> "2016-03-20 20:24:53"
#       (this also absorbs rows that used to read "2016-03-20 20:24:32")
#   - de-de sheet, column E ("Latest Handoff Datetime"):
#       "2016-03-20 20:24:12" -> "2016-03-20 20:24:57"
#       (this also absorbs rows that used to read "2016-03-20 20:24:35")
#
# Updating every row that shared the old timestamp (rather than only the
# rows that used the now-retired duplicate) keeps the shared-string table
# consistent: the stale duplicate strings become unreferenced and are
# dropped on save, while the surviving string's text is simply updated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: column D, rows 7,10,11,12,13,14,15,16
$overviewRows = 7,10,11,12,13,14,15,16
foreach ($r in $overviewRows) {
    $overview.Cells.Item($r, 4).Value2 = "2016-24-20 20:24:57"
}

# zh-cn sheet: column E, rows 7,10,11,12,13,14,15,16
$zhcnRows = 7,10,11,12,13,14,15,16
foreach ($r in $zhcnRows) {
    $zhcn.Cells.Item($r, 5).Value2 = "2016-03-20 20:24:53"
}

# de-de sheet: column E, rows 7,10,11,12,13,14,15,16
$dedeRows = 7,10,11,12,13,14,15,16
foreach ($r in $dedeRows) {
    $dede.Cells.Item($r, 5).Value2 = "2016-03-20 20:24:57"
}

Write-Output "Report regenerated for handoff."
